# Update figures ("想去人数" / "最低票价") on the "展览" (sheet 1) and
# "全部类型" (sheet 4) worksheets, matching the regenerated data snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(3, 6).Value = 5608   # F3 ("想去人数")
$ws.Cells.Item(3, 7).Value = 109   # G3 ("最低票价")
$ws.Cells.Item(5, 6).Value = 687   # F5 ("想去人数")
$ws.Cells.Item(6, 6).Value = 682   # F6 ("想去人数")
$ws.Cells.Item(7, 6).Value = 38   # F7 ("想去人数")
$ws.Cells.Item(9, 6).Value = 1086   # F9 ("想去人数")
$ws.Cells.Item(11, 6).Value = 1595   # F11 ("想去人数")
$ws.Cells.Item(12, 6).Value = 5669   # F12 ("想去人数")
$ws.Cells.Item(14, 6).Value = 303   # F14 ("想去人数")
$ws.Cells.Item(16, 6).Value = 65   # F16 ("想去人数")
$ws.Cells.Item(17, 6).Value = 33   # F17 ("想去人数")
$ws.Cells.Item(18, 6).Value = 119   # F18 ("想去人数")
$ws.Cells.Item(19, 6).Value = 4610   # F19 ("想去人数")
$ws.Cells.Item(20, 6).Value = 237   # F20 ("想去人数")
$ws.Cells.Item(21, 6).Value = 1220   # F21 ("想去人数")
$ws.Cells.Item(22, 6).Value = 125   # F22 ("想去人数")
$ws.Cells.Item(23, 6).Value = 84   # F23 ("想去人数")
$ws.Cells.Item(25, 6).Value = 80   # F25 ("想去人数")
$ws.Cells.Item(26, 6).Value = 219   # F26 ("想去人数")
$ws.Cells.Item(27, 6).Value = 76   # F27 ("想去人数")
$ws.Cells.Item(28, 6).Value = 156   # F28 ("想去人数")
$ws.Cells.Item(30, 6).Value = 358   # F30 ("想去人数")
$ws.Cells.Item(31, 6).Value = 44   # F31 ("想去人数")
$ws.Cells.Item(33, 6).Value = 71   # F33 ("想去人数")
$ws.Cells.Item(34, 6).Value = 22   # F34 ("想去人数")
$ws.Cells.Item(36, 6).Value = 7   # F36 ("想去人数")
$ws.Cells.Item(37, 6).Value = 48   # F37 ("想去人数")

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(3, 6).Value = 5608   # F3 ("想去人数")
$ws.Cells.Item(3, 7).Value = 109   # G3 ("最低票价")
$ws.Cells.Item(5, 6).Value = 687   # F5 ("想去人数")
$ws.Cells.Item(6, 6).Value = 682   # F6 ("想去人数")
$ws.Cells.Item(7, 6).Value = 38   # F7 ("想去人数")
$ws.Cells.Item(9, 6).Value = 1086   # F9 ("想去人数")
$ws.Cells.Item(11, 6).Value = 1595   # F11 ("想去人数")
$ws.Cells.Item(12, 6).Value = 5669   # F12 ("想去人数")
$ws.Cells.Item(14, 6).Value = 303   # F14 ("想去人数")
$ws.Cells.Item(16, 6).Value = 65   # F16 ("想去人数")
$ws.Cells.Item(17, 6).Value = 33   # F17 ("想去人数")
$ws.Cells.Item(18, 6).Value = 119   # F18 ("想去人数")
$ws.Cells.Item(19, 6).Value = 4610   # F19 ("想去人数")
$ws.Cells.Item(20, 6).Value = 237   # F20 ("想去人数")
$ws.Cells.Item(21, 6).Value = 1220   # F21 ("想去人数")
$ws.Cells.Item(22, 6).Value = 125   # F22 ("想去人数")
$ws.Cells.Item(23, 6).Value = 84   # F23 ("想去人数")
$ws.Cells.Item(25, 6).Value = 80   # F25 ("想去人数")
$ws.Cells.Item(26, 6).Value = 219   # F26 ("想去人数")
$ws.Cells.Item(27, 6).Value = 76   # F27 ("想去人数")
$ws.Cells.Item(28, 6).Value = 156   # F28 ("想去人数")
$ws.Cells.Item(30, 6).Value = 358   # F30 ("想去人数")
$ws.Cells.Item(31, 6).Value = 44   # F31 ("想去人数")
$ws.Cells.Item(33, 6).Value = 71   # F33 ("想去人数")
$ws.Cells.Item(34, 6).Value = 22   # F34 ("想去人数")
$ws.Cells.Item(36, 6).Value = 7   # F36 ("想去人数")
$ws.Cells.Item(37, 6).Value = 48   # F37 ("想去人数")

